$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date stamp in the header (A1), moved forward one month
# (45406 -> 45436, serialized as a date via style 23 / numFmtId 14).
$ws.Range("A1").Value = 45436

# Update the price list (column D) for rows 33-38.
$ws.Range("D33").Value = 4282.249
$ws.Range("D34").Value = 4282.249
$ws.Range("D35").Value = 5601.891
$ws.Range("D36").Value = 5601.932
$ws.Range("D37").Value = 5601.891
$ws.Range("D38").Value = 7061.332
